$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 228
$ws.Cells.Item(4, 6).Value = 481
$ws.Cells.Item(5, 6).Value = 2033
$ws.Cells.Item(7, 6).Value = 7788
$ws.Cells.Item(9, 6).Value = 30
$ws.Cells.Item(11, 6).Value = 226
$ws.Cells.Item(12, 6).Value = 1747
$ws.Cells.Item(16, 6).Value = 3846
$ws.Cells.Item(17, 6).Value = 5946
$ws.Cells.Item(18, 6).Value = 672
$ws.Cells.Item(20, 6).Value = 1061
$ws.Cells.Item(22, 6).Value = 403
$ws.Cells.Item(23, 6).Value = 6120
$ws.Cells.Item(26, 6).Value = 4164
$ws.Cells.Item(27, 6).Value = 696
$ws.Cells.Item(28, 6).Value = 1917
$ws.Cells.Item(30, 6).Value = 291
$ws.Cells.Item(32, 6).Value = 4
$ws.Cells.Item(33, 6).Value = 30
$ws.Cells.Item(35, 6).Value = 28
$ws.Cells.Item(38, 6).Value = 492
$ws.Cells.Item(41, 6).Value = 394
$ws.Cells.Item(43, 6).Value = 1108
$ws.Cells.Item(44, 6).Value = 552

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 987
$ws.Cells.Item(12, 6).Value = 354
$ws.Cells.Item(13, 6).Value = 397
$ws.Cells.Item(16, 6).Value = 105
$ws.Cells.Item(19, 6).Value = 344
$ws.Cells.Item(20, 6).Value = 161
$ws.Cells.Item(22, 6).Value = 52
$ws.Cells.Item(23, 6).Value = 21

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 446
$ws.Cells.Item(6, 6).Value = 1551
$ws.Cells.Item(7, 6).Value = 463
$ws.Cells.Item(9, 6).Value = 879
$ws.Cells.Item(10, 6).Value = 1038
$ws.Cells.Item(11, 6).Value = 1197
$ws.Cells.Item(12, 6).Value = 1534

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 1551
$ws.Cells.Item(3, 6).Value = 228
$ws.Cells.Item(4, 6).Value = 481
$ws.Cells.Item(5, 6).Value = 463
$ws.Cells.Item(7, 6).Value = 2033
$ws.Cells.Item(8, 6).Value = 7788
$ws.Cells.Item(10, 6).Value = 879
$ws.Cells.Item(11, 6).Value = 226
$ws.Cells.Item(12, 6).Value = 1747
$ws.Cells.Item(14, 6).Value = 1197
$ws.Cells.Item(18, 6).Value = 1534
$ws.Cells.Item(19, 6).Value = 3846
$ws.Cells.Item(20, 6).Value = 354
$ws.Cells.Item(21, 6).Value = 397
$ws.Cells.Item(22, 6).Value = 672
$ws.Cells.Item(23, 6).Value = 1061
$ws.Cells.Item(25, 6).Value = 403
$ws.Cells.Item(26, 6).Value = 6120
$ws.Cells.Item(28, 6).Value = 4164
$ws.Cells.Item(29, 6).Value = 696
$ws.Cells.Item(30, 6).Value = 1917
$ws.Cells.Item(32, 6).Value = 291
$ws.Cells.Item(33, 6).Value = 30
$ws.Cells.Item(34, 6).Value = 161
$ws.Cells.Item(38, 6).Value = 492
$ws.Cells.Item(40, 6).Value = 21
$ws.Cells.Item(42, 6).Value = 394
$ws.Cells.Item(43, 6).Value = 1108
$ws.Cells.Item(44, 6).Value = 552
